$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 113.666664
$ws.Range("I5").Value = 113.666664
$ws.Range("K5").Value = 113.666664
$ws.Range("M5").Value = 1.333336000000003

$ws.Range("H74").Value = 1620.6666
$ws.Range("I74").Value = 1620.6666
$ws.Range("K74").Value = 1620.6666
$ws.Range("M74").Value = -684.6666

$ws.Range("H77").Value = 1620.6666
$ws.Range("I77").Value = 1620.6666
$ws.Range("K77").Value = 8103.333000000001
$ws.Range("M77").Value = -3423.333000000001

$ws.Range("H112").Value = 1063.12
$ws.Range("J112").Value = 1063.12
$ws.Range("L112").Value = 3189.36
$ws.Range("N112").Value = -5405.36

$ws.Range("H129").Value = 900.56604
$ws.Range("I129").Value = 1021
$ws.Range("J129").Value = 882.23914
$ws.Range("K129").Value = 3063
$ws.Range("L129").Value = 2646.71742
$ws.Range("M129").Value = 1937
$ws.Range("N129").Value = -12646.71742

$ws.Range("H133").Value = 89000
$ws.Range("J133").Value = 89000
$ws.Range("L133").Value = 89000
$ws.Range("N133").Value = -99120

$ws.Range("H137").Value = 1760.5883
$ws.Range("I137").Value = 1107.4286
$ws.Range("J137").Value = 1929.9259
$ws.Range("K137").Value = 3322.2858
$ws.Range("L137").Value = 5789.7777
$ws.Range("M137").Value = -772.2857999999997
$ws.Range("N137").Value = -10889.7777

$ws.Range("H138").Value = 2076.5557
$ws.Range("I138").Value = 1653.2222
$ws.Range("J138").Value = 2499.889
$ws.Range("K138").Value = 4959.6666
$ws.Range("L138").Value = 7499.667
$ws.Range("M138").Value = 180.3334000000004
$ws.Range("N138").Value = -17779.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 19800
$ws.Range("J9").Value = 19800
$ws.Range("L9").Value = 19800
$ws.Range("N9").Value = -20140

$ws.Range("H20").Value = 19800
$ws.Range("J20").Value = 19800
$ws.Range("L20").Value = 19800
$ws.Range("N20").Value = -20340

$ws.Range("H32").Value = 4770.3447
$ws.Range("I32").Value = 3751.818
$ws.Range("J32").Value = 7971.4287
$ws.Range("K32").Value = 3751.818
$ws.Range("L32").Value = 7971.4287
$ws.Range("M32").Value = -3464.818
$ws.Range("N32").Value = -8545.4287

$ws.Range("H61").Value = 2821.2222
$ws.Range("I61").Value = 1854.7222
$ws.Range("K61").Value = 1854.7222
$ws.Range("M61").Value = -1642.7222

$ws.Range("H74").Value = 1300.9546
$ws.Range("I74").Value = 593.86664
$ws.Range("K74").Value = 593.86664
$ws.Range("M74").Value = 280.13336

$ws.Range("H77").Value = 1300.9546
$ws.Range("I77").Value = 593.86664
$ws.Range("K77").Value = 2969.3332
$ws.Range("M77").Value = 1398.6668

$ws.Range("H132").Value = 1937
$ws.Range("I132").Value = 1903.7812
$ws.Range("K132").Value = 5711.3436
$ws.Range("M132").Value = -3181.3436

$ws.Range("H136").Value = 2821.2222
$ws.Range("I136").Value = 1854.7222
$ws.Range("K136").Value = 5564.1666
$ws.Range("M136").Value = -3014.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 45936
$ws.Range("J76").Value = 45936
$ws.Range("L76").Value = 45936
$ws.Range("N76").Value = -46566

$ws.Range("H79").Value = 45936
$ws.Range("J79").Value = 45936
$ws.Range("L79").Value = 45936
$ws.Range("N79").Value = -48120

$ws.Range("H86").Value = 2000000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2000000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2000000
$ws.Range("M86").Value = $null
$ws.Range("N86").Value = -2002246

$ws.Range("H89").Value = 2000000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2000000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 10000000
$ws.Range("M89").Value = $null
$ws.Range("N89").Value = -10011232

$ws.Range("H132").Value = 54998.5
$ws.Range("J132").Value = 54998.5
$ws.Range("L132").Value = 54998.5
$ws.Range("N132").Value = -65118.5

$ws.Range("H134").Value = 5416.593
$ws.Range("I134").Value = 5529.92
$ws.Range("K134").Value = 16589.76
$ws.Range("M134").Value = -14054.76

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1394.5
$ws.Range("I31").Value = 980.3333
$ws.Range("J31").Value = 1490.0769
$ws.Range("K31").Value = 980.3333
$ws.Range("L31").Value = 1490.0769
$ws.Range("M31").Value = -685.3333
$ws.Range("N31").Value = -2080.0769

$ws.Range("H34").Value = 1394.5
$ws.Range("I34").Value = 980.3333
$ws.Range("J34").Value = 1490.0769
$ws.Range("K34").Value = 980.3333
$ws.Range("L34").Value = 1490.0769
$ws.Range("M34").Value = -778.3333
$ws.Range("N34").Value = -1894.0769

$ws.Range("H132").Value = 2060.5334
$ws.Range("I132").Value = 1810.9286
$ws.Range("K132").Value = 5432.7858
$ws.Range("M132").Value = -2902.7858

$ws.Range("H134").Value = 1331.3448
$ws.Range("I134").Value = 1048.4615
$ws.Range("K134").Value = 3145.3845
$ws.Range("M134").Value = -610.3844999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 2999.6667
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 2999.6667
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 8999.000100000001
$ws.Range("M116").Value = $null
$ws.Range("N116").Value = -15883.0001

$ws.Range("H129").Value = 45685.125
$ws.Range("I129").Value = 595.7778
$ws.Range("J129").Value = 103657.14
$ws.Range("K129").Value = 1787.3334
$ws.Range("L129").Value = 310971.42
$ws.Range("M129").Value = 3212.6666
$ws.Range("N129").Value = -320971.42

$ws.Range("H131").Value = 11925093
$ws.Range("J131").Value = 23031.055
$ws.Range("L131").Value = 69093.16500000001
$ws.Range("N131").Value = -79173.16500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -56884

$ws.Range("H132").Value = 1203677.1
$ws.Range("I132").Value = 1540014.4
$ws.Range("K132").Value = 4620043.199999999
$ws.Range("M132").Value = -4617513.199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = $null

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = $null

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = $null

$ws.Range("H122").Value = 8752.846
$ws.Range("I122").Value = 6911.25
$ws.Range("J122").Value = 11699.4
$ws.Range("K122").Value = 20733.75
$ws.Range("L122").Value = 35098.2
$ws.Range("M122").Value = -18283.75
$ws.Range("N122").Value = -39998.2

$ws.Range("H132").Value = 3795
$ws.Range("I132").Value = 1298.5
$ws.Range("K132").Value = 3895.5
$ws.Range("M132").Value = -1365.5

$ws.Range("H136").Value = 3402
$ws.Range("I136").Value = 2182.7827
$ws.Range("J136").Value = 6517.778
$ws.Range("K136").Value = 6548.348100000001
$ws.Range("L136").Value = 19553.334
$ws.Range("M136").Value = -3998.348100000001
$ws.Range("N136").Value = -24653.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 13676
$ws.Range("I126").Value = 19411.857
$ws.Range("J126").Value = 6984.1665
$ws.Range("K126").Value = 58235.571
$ws.Range("L126").Value = 20952.4995
$ws.Range("M126").Value = -55765.571
$ws.Range("N126").Value = -25892.4995

$ws.Range("H132").Value = 1474.225
$ws.Range("I132").Value = 1252.4
$ws.Range("J132").Value = 2139.7
$ws.Range("K132").Value = 3757.2
$ws.Range("L132").Value = 6419.099999999999
$ws.Range("M132").Value = -1227.2
$ws.Range("N132").Value = -11479.1

$ws.Range("H136").Value = 11577256
$ws.Range("I136").Value = 19160034
$ws.Range("J136").Value = 3541.3157
$ws.Range("K136").Value = 57480102
$ws.Range("L136").Value = 10623.9471
$ws.Range("M136").Value = -57477552
$ws.Range("N136").Value = -15723.9471
